$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4115.3076
$ws.Range("I40").Value = 1524.75
$ws.Range("J40").Value = 5266.6665
$ws.Range("K40").Value = 1524.75
$ws.Range("L40").Value = 5266.6665
$ws.Range("M40").Value = -1349.75
$ws.Range("N40").Value = -5616.6665
$ws.Range("H113").Value = 18455.092
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H116").Value = 510470.84
$ws.Range("I116").Value = 911782.25
$ws.Range("J116").Value = 19979.111
$ws.Range("K116").Value = 911782.25
$ws.Range("L116").Value = 19979.111
$ws.Range("M116").Value = -908340.25
$ws.Range("N116").Value = -26863.111
$ws.Range("H137").Value = 2275.578
$ws.Range("I137").Value = 1556.7838
$ws.Range("J137").Value = 5600
$ws.Range("K137").Value = 4670.3514
$ws.Range("L137").Value = 16800
$ws.Range("M137").Value = -2120.3514
$ws.Range("N137").Value = -21900
$ws.Range("H138").Value = 2609.9111
$ws.Range("I138").Value = 1525.409
$ws.Range("J138").Value = 2960.7793
$ws.Range("K138").Value = 4576.227000000001
$ws.Range("L138").Value = 8882.3379
$ws.Range("M138").Value = 563.7729999999992
$ws.Range("N138").Value = -19162.3379

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H45").Value = 2560
$ws.Range("I45").Value = 1200
$ws.Range("J45").Value = 2696
$ws.Range("K45").Value = 1200
$ws.Range("L45").Value = 2696
$ws.Range("M45").Value = -823
$ws.Range("N45").Value = -3450
$ws.Range("H102").Value = 1367.8334
$ws.Range("I102").Value = 1367.8334
$ws.Range("K102").Value = 1367.8334
$ws.Range("M102").Value = 254.1666
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 3006.5715
$ws.Range("I122").Value = 1701.7142
$ws.Range("J122").Value = 4311.4287
$ws.Range("K122").Value = 5105.142599999999
$ws.Range("L122").Value = 12934.2861
$ws.Range("M122").Value = -2655.142599999999
$ws.Range("N122").Value = -17834.2861
$ws.Range("H132").Value = 2057.577
$ws.Range("I132").Value = 992.0625
$ws.Range("J132").Value = 3762.4
$ws.Range("K132").Value = 2976.1875
$ws.Range("L132").Value = 11287.2
$ws.Range("M132").Value = -446.1875
$ws.Range("N132").Value = -16347.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H134").Value = 3173.6978
$ws.Range("I134").Value = 1813.9706
$ws.Range("J134").Value = 8310.444
$ws.Range("K134").Value = 5441.9118
$ws.Range("L134").Value = 24931.332
$ws.Range("M134").Value = -2906.9118
$ws.Range("N134").Value = -30001.332

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5052438.5
$ws.Range("I16").Value = 11112044
$ws.Range("K16").Value = 11112044
$ws.Range("M16").Value = -11111757
$ws.Range("H21").Value = 39999.5
$ws.Range("J21").Value = 39999.5
$ws.Range("L21").Value = 39999.5
$ws.Range("N21").Value = -40469.5
$ws.Range("H35").Value = 34999.25
$ws.Range("I35").Value = 4998
$ws.Range("J35").Value = 44999.668
$ws.Range("K35").Value = 4998
$ws.Range("L35").Value = 44999.668
$ws.Range("M35").Value = -4704
$ws.Range("N35").Value = -45587.668
$ws.Range("H36").Value = 18504.8
$ws.Range("I36").Value = 11262
$ws.Range("J36").Value = 23333.334
$ws.Range("K36").Value = 11262
$ws.Range("L36").Value = 23333.334
$ws.Range("M36").Value = -10874
$ws.Range("N36").Value = -24109.334
$ws.Range("H40").Value = 18504.8
$ws.Range("I40").Value = 11262
$ws.Range("J40").Value = 23333.334
$ws.Range("K40").Value = 11262
$ws.Range("L40").Value = 23333.334
$ws.Range("M40").Value = -11102
$ws.Range("N40").Value = -23653.334
$ws.Range("H41").Value = 34532.375
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 34532.375
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 34532.375
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -35388.375
$ws.Range("H42").Value = 50000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 50000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 50000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -51186
$ws.Range("H56").Value = 30103
$ws.Range("J56").Value = 30103
$ws.Range("L56").Value = 30103
$ws.Range("N56").Value = -31793
$ws.Range("H99").Value = 5098.25
$ws.Range("I99").Value = 3128.8572
$ws.Range("J99").Value = 6630
$ws.Range("K99").Value = 3128.8572
$ws.Range("L99").Value = 6630
$ws.Range("M99").Value = -1630.8572
$ws.Range("N99").Value = -9626
$ws.Range("H113").Value = 5052438.5
$ws.Range("I113").Value = 11112044
$ws.Range("K113").Value = 11112044
$ws.Range("M113").Value = -11109874
$ws.Range("H126").Value = 5098.25
$ws.Range("I126").Value = 3128.8572
$ws.Range("J126").Value = 6630
$ws.Range("K126").Value = 9386.571599999999
$ws.Range("L126").Value = 19890
$ws.Range("M126").Value = -6916.571599999999
$ws.Range("N126").Value = -24830
$ws.Range("H134").Value = 5970.6924
$ws.Range("I134").Value = 6364.421
$ws.Range("J134").Value = 4902
$ws.Range("K134").Value = 19093.263
$ws.Range("L134").Value = 14706
$ws.Range("M134").Value = -16558.263
$ws.Range("N134").Value = -19776

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1524.0869
$ws.Range("I5").Value = 662.0645
$ws.Range("J5").Value = 3305.6
$ws.Range("K5").Value = 1986.1935
$ws.Range("L5").Value = 9916.799999999999
$ws.Range("M5").Value = -1874.1935
$ws.Range("N5").Value = -10140.8
$ws.Range("H39").Value = 8592.691999999999
$ws.Range("I39").Value = 1501
$ws.Range("J39").Value = 9183.666999999999
$ws.Range("K39").Value = 4503
$ws.Range("L39").Value = 27551.001
$ws.Range("M39").Value = -4209
$ws.Range("N39").Value = -28139.001
$ws.Range("H55").Value = 4245
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 4993.3335
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 14980.0005
$ws.Range("M55").Value = -5823
$ws.Range("N55").Value = -15334.0005
$ws.Range("H131").Value = 8772941
$ws.Range("I131").Value = 41668012
$ws.Range("J131").Value = 922.04443
$ws.Range("K131").Value = 125004036
$ws.Range("L131").Value = 2766.13329
$ws.Range("M131").Value = -124998996
$ws.Range("N131").Value = -12846.13329
$ws.Range("H132").Value = 2918.0454
$ws.Range("I132").Value = 1127
$ws.Range("J132").Value = 3753.8667
$ws.Range("K132").Value = 10143
$ws.Range("L132").Value = 33784.8003
$ws.Range("M132").Value = -7613
$ws.Range("N132").Value = -38844.8003
$ws.Range("H135").Value = 1524.0869
$ws.Range("I135").Value = 662.0645
$ws.Range("J135").Value = 3305.6
$ws.Range("K135").Value = 5958.5805
$ws.Range("L135").Value = 29750.4
$ws.Range("M135").Value = -3423.5805
$ws.Range("N135").Value = -34820.39999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2734.1875
$ws.Range("I132").Value = 1669.5
$ws.Range("K132").Value = 5008.5
$ws.Range("M132").Value = -2478.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 75395.22
$ws.Range("J46").Value = 75395.22
$ws.Range("L46").Value = 75395.22
$ws.Range("N46").Value = -75857.22
$ws.Range("H113").Value = 8111.385
$ws.Range("I113").Value = 14689.429
$ws.Range("J113").Value = 437
$ws.Range("K113").Value = 44068.287
$ws.Range("L113").Value = 1311
$ws.Range("M113").Value = -41898.287
$ws.Range("N113").Value = -5651
$ws.Range("H123").Value = 29940
$ws.Range("J123").Value = 29940
$ws.Range("L123").Value = 29940
$ws.Range("N123").Value = -39740
$ws.Range("H126").Value = 7433.033
$ws.Range("I126").Value = 3591.8235
$ws.Range("J126").Value = 12456.154
$ws.Range("K126").Value = 10775.4705
$ws.Range("L126").Value = 37368.462
$ws.Range("M126").Value = -8305.470499999999
$ws.Range("N126").Value = -42308.462
$ws.Range("H131").Value = 42715
$ws.Range("J131").Value = 42715
$ws.Range("L131").Value = 42715
$ws.Range("N131").Value = -52795
$ws.Range("H132").Value = 7411374.5
$ws.Range("I132").Value = 4672.52
$ws.Range("J132").Value = 16669752
$ws.Range("K132").Value = 14017.56
$ws.Range("L132").Value = 50009256
$ws.Range("M132").Value = -11487.56
$ws.Range("N132").Value = -50014316
$ws.Range("H134").Value = 75395.22
$ws.Range("J134").Value = 75395.22
$ws.Range("L134").Value = 226185.66
$ws.Range("N134").Value = -231255.66
$ws.Range("H136").Value = 6222.353
$ws.Range("I136").Value = 2382.8333
$ws.Range("J136").Value = 8316.637000000001
$ws.Range("K136").Value = 7148.499899999999
$ws.Range("L136").Value = 24949.911
$ws.Range("M136").Value = -4598.499899999999
$ws.Range("N136").Value = -30049.911
